$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: drop the "Time" label from A1 (headers in B1/C1/D1 stay the same) ---
$ws.Range("A1").ClearContents()

# --- Fix the column-A timestamps: the original file accumulated float drift
#     (e.g. 43935.08333321759); recompute the clean "day + hour/24" value for
#     every existing row (2-289) and prep the format for the rows being added. ---
$ws.Range("A290:A313").NumberFormat = "m/d/yy h:mm"

$baseDay = 43935
for ($row = 2; $row -le 313; $row++) {
    $idx = $row - 2
    $dayOffset = [math]::Floor($idx / 24)
    $hourOffset = $idx % 24
    $val = $baseDay + $dayOffset + ($hourOffset / 24.0)
    $ws.Cells.Item($row, 1).Value = $val
}

# --- Append the new auction data for 2020-04-13 (rows 290-313) ---
$newData = @(
    @(18.02,22.3,16.23),
    @(16.32,19.79,15.76),
    @(15.26,20.76,14.55),
    @(15.21,21.06,12.42),
    @(15.06,22.23,11.92),
    @(15.2,21.95,11.8),
    @(15.2,18.81,12),
    @(14.89,18.09,11.99),
    @(14.42,16.06,12.12),
    @(14.51,15.01,12.42),
    @(11.94,13,13),
    @(12.03,14.38,14.38),
    @(10.59,15,15),
    @(9.01,11.9,11.9),
    @(5.81,5.56,5.56),
    @(5.61,6.93,6.93),
    @(7.9,12.42,12.42),
    @(15.13,16.15,12.42),
    @(23.73,24.82,14),
    @(23.97,27.97,17.07),
    @(24.83,31.84,19.46),
    @(24.96,27.49,24.1),
    @(25.66,26.06,25.15),
    @(23.59,23.76,23.76)
)

$startRow = 290
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $vals = $newData[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
}

# --- Match the saved selection / scroll state ---
$ws.Range("E18").Select()
